$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 12).Value = 3221
$ws.Cells.Item(3, 12).Value = 3322
$ws.Cells.Item(4, 12).Value = 838
$ws.Cells.Item(5, 12).Value = 189
$ws.Cells.Item(6, 12).Value = 2936
$ws.Cells.Item(7, 12).Value = 10506

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(6, 12).Value = 41
$ws.Cells.Item(7, 12).Value = 120

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 12).Value = 196
$ws.Cells.Item(3, 12).Value = 218
$ws.Cells.Item(6, 12).Value = 189
$ws.Cells.Item(7, 12).Value = 670

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(2, 12).Value = 79
$ws.Cells.Item(3, 12).Value = 100
$ws.Cells.Item(7, 12).Value = 247

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 12).Value = 137
$ws.Cells.Item(3, 12).Value = 147
$ws.Cells.Item(4, 12).Value = 25
$ws.Cells.Item(6, 12).Value = 167
$ws.Cells.Item(7, 12).Value = 484

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(2, 12).Value = 57
$ws.Cells.Item(6, 12).Value = 30
$ws.Cells.Item(7, 12).Value = 145

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(3, 12).Value = 113
$ws.Cells.Item(6, 12).Value = 116
$ws.Cells.Item(7, 12).Value = 380

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(6, 12).Value = 59
$ws.Cells.Item(7, 12).Value = 203

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(5, 12).Value = 3
$ws.Cells.Item(7, 12).Value = 179

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Cells.Item(6, 12).Value = 25
$ws.Cells.Item(7, 12).Value = 50

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(6, 12).Value = 83
$ws.Cells.Item(7, 12).Value = 350
$ws.Cells.Item(8, 12).Value = 670
$ws.Cells.Item(11, 12).Value = 176
$ws.Cells.Item(15, 12).Value = 79
$ws.Cells.Item(16, 12).Value = 23
$ws.Cells.Item(18, 12).Value = 79
$ws.Cells.Item(19, 12).Value = 297
$ws.Cells.Item(20, 12).Value = 265
$ws.Cells.Item(21, 12).Value = 31
$ws.Cells.Item(22, 12).Value = 32
$ws.Cells.Item(23, 12).Value = 112
$ws.Cells.Item(24, 12).Value = 25
$ws.Cells.Item(25, 12).Value = 55
$ws.Cells.Item(27, 12).Value = 101
$ws.Cells.Item(29, 12).Value = 568
$ws.Cells.Item(30, 12).Value = 50
$ws.Cells.Item(33, 12).Value = 484
$ws.Cells.Item(34, 12).Value = 68
$ws.Cells.Item(37, 12).Value = 380
$ws.Cells.Item(40, 12).Value = 27
$ws.Cells.Item(42, 12).Value = 335
$ws.Cells.Item(43, 12).Value = 81
$ws.Cells.Item(46, 12).Value = 23
$ws.Cells.Item(48, 12).Value = 143
$ws.Cells.Item(52, 12).Value = 210
$ws.Cells.Item(53, 12).Value = 120
$ws.Cells.Item(54, 12).Value = 219
$ws.Cells.Item(57, 12).Value = 37
$ws.Cells.Item(60, 12).Value = 63
$ws.Cells.Item(63, 12).Value = 33
$ws.Cells.Item(64, 12).Value = 70
$ws.Cells.Item(65, 12).Value = 203
$ws.Cells.Item(67, 12).Value = 386
$ws.Cells.Item(72, 12).Value = 51
$ws.Cells.Item(73, 12).Value = 90
$ws.Cells.Item(78, 12).Value = 131
$ws.Cells.Item(79, 12).Value = 269
$ws.Cells.Item(83, 12).Value = 247
$ws.Cells.Item(84, 12).Value = 104
$ws.Cells.Item(85, 12).Value = 532
$ws.Cells.Item(87, 12).Value = 36
$ws.Cells.Item(88, 12).Value = 124
$ws.Cells.Item(91, 12).Value = 153
$ws.Cells.Item(93, 12).Value = 58
$ws.Cells.Item(94, 12).Value = 125
$ws.Cells.Item(95, 12).Value = 145
$ws.Cells.Item(99, 12).Value = 179
$ws.Cells.Item(101, 12).Value = 10506

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(2, 12).Value = 114
$ws.Cells.Item(3, 12).Value = 146
$ws.Cells.Item(7, 12).Value = 386

$ws = $wb.Worksheets.Item('South Deering')
$ws.Cells.Item(6, 12).Value = 23
$ws.Cells.Item(7, 12).Value = 104

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(2, 12).Value = 47
$ws.Cells.Item(7, 12).Value = 219

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(3, 12).Value = 215
$ws.Cells.Item(7, 12).Value = 568

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(2, 12).Value = 19
$ws.Cells.Item(3, 12).Value = 33
$ws.Cells.Item(7, 12).Value = 143

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(3, 12).Value = 89
$ws.Cells.Item(7, 12).Value = 297

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Cells.Item(2, 12).Value = 34
$ws.Cells.Item(3, 12).Value = 25
$ws.Cells.Item(7, 12).Value = 83

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(3, 12).Value = 104
$ws.Cells.Item(6, 12).Value = 94
$ws.Cells.Item(7, 12).Value = 335

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(2, 12).Value = 38
$ws.Cells.Item(3, 12).Value = 39
$ws.Cells.Item(6, 12).Value = 38
$ws.Cells.Item(7, 12).Value = 131

$ws = $wb.Worksheets.Item('Dunning')
$ws.Cells.Item(3, 12).Value = 6
$ws.Cells.Item(7, 12).Value = 25

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Cells.Item(3, 12).Value = 8
$ws.Cells.Item(7, 12).Value = 23

$ws = $wb.Worksheets.Item('Douglas')
$ws.Cells.Item(3, 12).Value = 43
$ws.Cells.Item(6, 12).Value = 26
$ws.Cells.Item(7, 12).Value = 112

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Cells.Item(3, 12).Value = 63
$ws.Cells.Item(7, 12).Value = 153

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Cells.Item(4, 12).Value = 3
$ws.Cells.Item(7, 12).Value = 31

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(2, 12).Value = 92
$ws.Cells.Item(3, 12).Value = 95
$ws.Cells.Item(7, 12).Value = 269

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Cells.Item(6, 12).Value = 20
$ws.Cells.Item(7, 12).Value = 70

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(2, 12).Value = 85
$ws.Cells.Item(3, 12).Value = 80
$ws.Cells.Item(4, 12).Value = 25
$ws.Cells.Item(7, 12).Value = 265

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Cells.Item(2, 12).Value = 30
$ws.Cells.Item(3, 12).Value = 30
$ws.Cells.Item(5, 12).Value = 2
$ws.Cells.Item(6, 12).Value = 10
$ws.Cells.Item(7, 12).Value = 79

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Cells.Item(3, 12).Value = 15
$ws.Cells.Item(7, 12).Value = 58

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(6, 12).Value = 99
$ws.Cells.Item(7, 12).Value = 350

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Cells.Item(4, 12).Value = 7
$ws.Cells.Item(6, 12).Value = 25
$ws.Cells.Item(7, 12).Value = 68

$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(6, 12).Value = 44
$ws.Cells.Item(7, 12).Value = 125

$ws = $wb.Worksheets.Item('East Side')
$ws.Cells.Item(3, 12).Value = 29
$ws.Cells.Item(7, 12).Value = 55

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(6, 12).Value = 17
$ws.Cells.Item(7, 12).Value = 79

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(2, 12).Value = 65
$ws.Cells.Item(3, 12).Value = 55
$ws.Cells.Item(6, 12).Value = 41
$ws.Cells.Item(7, 12).Value = 176

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(6, 12).Value = 22
$ws.Cells.Item(7, 12).Value = 90

$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(3, 12).Value = 45
$ws.Cells.Item(7, 12).Value = 124

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(3, 12).Value = 34
$ws.Cells.Item(7, 12).Value = 101

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Cells.Item(6, 12).Value = 12
$ws.Cells.Item(7, 12).Value = 37

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Cells.Item(3, 12).Value = 24
$ws.Cells.Item(7, 12).Value = 63

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(4, 12).Value = 15
$ws.Cells.Item(7, 12).Value = 81

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(6, 12).Value = 111
$ws.Cells.Item(7, 12).Value = 532

$ws = $wb.Worksheets.Item('Clearing')
$ws.Cells.Item(3, 12).Value = 14
$ws.Cells.Item(7, 12).Value = 32

$ws = $wb.Worksheets.Item('Old Town')
$ws.Cells.Item(3, 12).Value = 12
$ws.Cells.Item(7, 12).Value = 51

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Cells.Item(2, 12).Value = 7
$ws.Cells.Item(7, 12).Value = 27

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(2, 12).Value = 74
$ws.Cells.Item(7, 12).Value = 210

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Cells.Item(4, 12).Value = 5
$ws.Cells.Item(7, 12).Value = 36

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(7, 12).Value = 23
